$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing checkmark emoji from names, and fix a name swap in rows 12-13
$ws.Range("C2").Value = "Jonnah"
$ws.Range("C3").Value = "Mandy"
$ws.Range("C4").Value = "Jonnah"
$ws.Range("C5").Value = "Sam"
$ws.Range("C7").Value = "Sungwoo"
$ws.Range("C8").Value = "Minal"
$ws.Range("C9").Value = "yujin"
$ws.Range("C11").Value = "Mandy"
$ws.Range("C12").Value = "Minal"
$ws.Range("C13").Value = "Fionna"
$ws.Range("C15").Value = "yujin"
$ws.Range("C16").Value = "Sam"
$ws.Range("C17").Value = "Sungwoo"
$ws.Range("C18").Value = "Fionna"
